$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Title (appears twice: Heading1 at top, bold run near the bottom)
Replace-Text "Play Fortune Caravan Free: Review and Gameplay" "Play Fortune Caravan for Free"

# "What we like" bullets
Replace-Text "Intuitive interface making it accessible for all kinds of players" "Intuitive interface for all players"
Replace-Text "Mobile versions make it possible to play on the go" "Gamble feature and Scatter icons add excitement"
Replace-Text "Two Scatter icons make the game more fun and engaging" "Wide betting range suitable for different players"
Replace-Text "Favorable Return to Player Percentage (RTP) of 96.30%" "Compatible with desktop and mobile devices"

# "What we don't like" bullets
Replace-Text "High-volatility slot game may not be suitable for some players" "High-volatility may not be ideal for all players"
Replace-Text "Wild symbol only appears on reel 5" "Limited availability of Wild symbol"

# Closing meta description (italic run)
Replace-Text "Read our Fortune Caravan review, featuring gameplay and features, symbol descriptions, betting range, compatibility, accessibility, and RTP. Play it for free now!" "Read our review of Fortune Caravan, a slot game with exciting gameplay and features. Play now for free."
